$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'AU-4 (1),AU-4'
$ws.Range("A3").Value = 'AU-14 (1),AU-4'
$ws.Range("A5").Value = 'SC-5 (2),SC-5,CM-6 b'
$ws.Range("A6").Value = 'AC-6 (8),CM-5 (1),AU-7 b,AU-12 (3),AU-7 a,AU-8 b,AC-6 (9)'
$ws.Range("A7").Value = 'CM-5 (1),AU-7 b,AU-12 (3),AU-12 c,AU-7 a,AU-8 b,AU-12 a,CM-6 b'
$ws.Range("A45").Value = 'IA-8,IA-2,AU-3 (1)'
$ws.Range("A46").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A47").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A48").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A49").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A50").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A51").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A52").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A53").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A54").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A55").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A56").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A57").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A58").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A59").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A60").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A61").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A62").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A63").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A64").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A65").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A66").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A67").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A68").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A69").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A70").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A71").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A72").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A73").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A74").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A75").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A76").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A77").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A78").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A79").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A80").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A81").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A82").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A83").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A84").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A85").Value = 'AU-3,MA-4 (1) (a),AU-3 (1)'
$ws.Range("A86").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A87").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A88").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A89").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A90").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A91").Value = 'AU-12 c,MA-4 (1) (a),AU-3 (1)'
$ws.Range("A92").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A93").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A94").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A95").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A96").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A97").Value = 'AU-12 c,AU-3,MA-4 (1) (a),AU-3 (1)'
$ws.Range("A98").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A99").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A100").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A101").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A102").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A103").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A104").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A105").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A106").Value = 'AU-12 c,AU-3,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A107").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A108").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A109").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A110").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A111").Value = 'AU-12 c,AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AC-2 (4)'
$ws.Range("A112").Value = 'AU-12 c,AU-14 (1),AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A113").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A114").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A116").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A120").Value = 'AU-12 c,AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A121").Value = 'AU-12 c,AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A127").Value = 'CM-5 (1),AU-12 c,AC-2 (4),AC-6 (9)'
$ws.Range("A129").Value = 'IA-5 (1) (b),IA-5 (1) (a),CM-6 b'
$ws.Range("A133").Value = 'SC-13,AC-17 (2),MA-4 c,SC-8'
$ws.Range("A134").Value = 'MA-4 (7),AC-12,SC-10,MA-4 e'
$ws.Range("A135").Value = 'AC-12,SC-10'
$ws.Range("A136").Value = 'AC-12,SC-10'
$ws.Range("A137").Value = 'AC-11 a,SC-10'
$ws.Range("A138").Value = 'CM-5 (1),AU-7 (1),AU-14 (1),AU-7 a,AU-3,AU-12 a,AU-3 (1),CM-6 b,AU-6 (4),MA-4 (1) (a)'
$ws.Range("A141").Value = 'AU-9 (3),AU-9'
$ws.Range("A142").Value = 'AU-9 (3),AU-9'
$ws.Range("A180").Value = 'CM-7 b,AC-17 (1),AC-17 (9),CM-6 b'
$ws.Range("A181").Value = 'CM-7 b,AC-17 (1),CM-6 b'
$ws.Range("A182").Value = 'AU-9,SI-11 b'
$ws.Range("A183").Value = 'AU-9,SI-11 b'
$ws.Range("A184").Value = 'AU-9,SI-11 b'
$ws.Range("A185").Value = 'AU-9,SI-11 b'
$ws.Range("A186").Value = 'AU-9,SI-11 b'
$ws.Range("A187").Value = 'AU-9,SI-11 b'
$ws.Range("A194").Value = 'AU-3,CM-6 b'
$ws.Range("A212").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A213").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A214").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A215").Value = 'AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Range("A222").Value = 'IA-2 (3),IA-2 (4),IA-2,IA-2 (2),IA-2 (5)'
$ws.Range("A223").Value = 'IA-2 (3),IA-2 (4),IA-2,IA-2 (2),IA-2 (5)'
$ws.Range("A224").Value = 'SC-8 (1),SC-8,AC-18 (1)'
$ws.Range("A226").Value = 'IA-5 (1) (c),IA-7'
$ws.Range("A227").Value = 'IA-7,CM-6 b'
$ws.Range("A228").Value = 'IA-7,CM-6 b'
$ws.Range("A229").Value = 'IA-7,CM-6 b'
$ws.Range("A231").Value = 'CM-7 a,IA-7'
$ws.Range("A232").Value = 'SC-13,MA-4 (6)'
$ws.Range("A234").Value = 'SC-13,MA-4 (6)'
$ws.Range("A242").Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Range("A243").Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Range("A245").Value = 'CM-6 b,SI-16,SC-2'
$ws.Range("A247").Value = 'SC-3,SI-16'
$ws.Range("A259").Value = 'IA-3,CM-6 b'
$ws.Range("A260").Value = 'IA-3,CM-6 b'
$ws.Range("A261").Value = 'IA-3,CM-6 b'
$ws.Range("A262").Value = 'IA-3,CM-6 b'
$ws.Range("A271").Value = 'IA-2 (4),IA-2 (2),IA-2 (1),IA-2 (3)'
$ws.Range("A273").Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Range("A276").Value = 'SC-4,CM-6 b'
$ws.Range("A277").Value = 'SC-4,SC-2'
$ws.Range("A278").Value = 'SC-4,SC-2'
$ws.Range("A281").Value = 'AU-12 a,CM-6 b'
$ws.Range("A284").Value = 'CM-5 (3),CM-6 b'
$ws.Range("A300").Value = 'IA-2 (11),IA-2 (1),IA-2 (12)'
$ws.Range("A310").Value = 'AU-8 (1) (b),AU-8 b,AU-8 (1) (a)'
$ws.Range("A328").Value = 'CM-5 (1),AU-12 c'
$ws.Range("A330").Value = 'AU-5 b,AU-5 a'
$ws.Range("A342").Value = 'CM-7 b,IA-3'
$ws.Range("A345").Value = 'CM-7 b,AC-17 (1)'
$ws.Range("A346").Value = 'CM-7 a,AC-18 (1)'
$ws.Range("A347").Value = 'IA-5 (1) (c),CM-7 a,CM-6 b'
$ws.Range("A361").Value = 'SI-6 d,SI-6 b,CM-3 (5)'
$ws.Range("A362").Value = 'CM-7 a,CM-6 b'
$ws.Range("A367").Value = 'CM-7 a,SI-16'
$ws.Range("A374").Value = 'CM-7 a,CM-6 b'
$ws.Range("A375").Value = 'CM-7 a,CM-6 b'
$ws.Range("A376").Value = 'CM-7 a,CM-6 b'
$ws.Range("A389").Value = 'SC-3,SI-6 a'
$ws.Range("A391").Value = 'IA-5 (1) (a),CM-6 b'
$ws.Range("A401").Value = 'SC-3,CM-6 b'
$ws.Range("A402").Value = 'SC-3,CM-6 b'
$ws.Range("A403").Value = 'SC-3,CM-6 b'
$ws.Range("A448").Value = 'IA-5 (1) (c),CM-6 b'
$ws.Range("A450").Value = 'CM-5 (1),CM-6 b'
$ws.Range("A451").Value = 'CM-5 (1),CM-6 b'
$ws.Range("A524").Value = 'CM-6 b,SC-2'
$ws.Range("A525").Value = 'CM-6 b,SC-2'
$ws.Range("A558").Value = 'CM-3 (5),SI-6 a'

# Fix wording in K361: "periodicly scan" -> "scan periodically"
$k361 = $ws.Range("K361").Value2
$k361 = $k361.Replace("If AIDE is not configured to periodicly scan then this is a finding.", "If AIDE is not configured to scan periodically then this is a finding.")
$ws.Range("K361").Value2 = $k361

